$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 743.6
$ws.Range("I19").Value = 704.2857
$ws.Range("J19").Value = 764.7692
$ws.Range("K19").Value = 704.2857
$ws.Range("L19").Value = 764.7692
$ws.Range("M19").Value = -529.2857
$ws.Range("N19").Value = -1114.7692

$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("M47").Value = -19028

$ws.Range("H48").Value = 2333.3333
$ws.Range("J48").Value = 3000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9584

$ws.Range("H56").Value = 2333.3333
$ws.Range("J56").Value = 3000
$ws.Range("L56").Value = 9000
$ws.Range("N56").Value = -10068

$ws.Range("H64").Value = 5813.8887
$ws.Range("I64").Value = 3280
$ws.Range("J64").Value = 6537.857
$ws.Range("K64").Value = 3280
$ws.Range("L64").Value = 6537.857
$ws.Range("N64").Value = -7033.857
$ws.Range("M64").Value = -3032

$ws.Range("H67").Value = 5813.8887
$ws.Range("I67").Value = 3280
$ws.Range("J67").Value = 6537.857
$ws.Range("K67").Value = 3280
$ws.Range("L67").Value = 6537.857
$ws.Range("N67").Value = -8253.857
$ws.Range("M67").Value = -2422

$ws.Range("H70").Value = 1802.7
$ws.Range("I70").Value = 1323.5
$ws.Range("J70").Value = 2122.1667
$ws.Range("K70").Value = 3970.5
$ws.Range("L70").Value = 6366.500100000001
$ws.Range("M70").Value = -3700.5
$ws.Range("N70").Value = -6906.500100000001

$ws.Range("H73").Value = 1802.7
$ws.Range("I73").Value = 1323.5
$ws.Range("J73").Value = 2122.1667
$ws.Range("K73").Value = 3970.5
$ws.Range("L73").Value = 6366.500100000001
$ws.Range("M73").Value = -3034.5
$ws.Range("N73").Value = -8238.500100000001

$ws.Range("H129").Value = 986.94543
$ws.Range("J129").Value = 1038.74
$ws.Range("L129").Value = 3116.22
$ws.Range("N129").Value = -13116.22

$ws.Range("H132").Value = 239567.69
$ws.Range("I132").Value = 258798.17
$ws.Range("J132").Value = 58801.2
$ws.Range("K132").Value = 776394.51
$ws.Range("L132").Value = 176403.6
$ws.Range("M132").Value = -773864.51
$ws.Range("N132").Value = -181463.6

$ws.Range("H137").Value = 2111.4285
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2111.4285
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 6334.2855
$ws.Range("N137").Value = -11434.2855
$ws.Range("M137").ClearContents()

$ws.Range("H140").Value = 30500
$ws.Range("J140").Value = 41000
$ws.Range("L140").Value = 41000
$ws.Range("N140").Value = -51360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2380.3281
$ws.Range("I32").Value = 1962.8536
$ws.Range("K32").Value = 1962.8536
$ws.Range("M32").Value = -1675.8536

$ws.Range("H76").Value = 29500.5
$ws.Range("J76").Value = 50001
$ws.Range("L76").Value = 50001
$ws.Range("N76").Value = -50677

$ws.Range("H79").Value = 29500.5
$ws.Range("J79").Value = 50001
$ws.Range("L79").Value = 50001
$ws.Range("N79").Value = -52341

$ws.Range("H122").Value = 2014.45
$ws.Range("I122").Value = 1775
$ws.Range("J122").Value = 2307.111
$ws.Range("K122").Value = 5325
$ws.Range("L122").Value = 6921.333
$ws.Range("M122").Value = -2875
$ws.Range("N122").Value = -11821.333

$ws.Range("H132").Value = 2382.3242
$ws.Range("I132").Value = 1930.4073
$ws.Range("J132").Value = 3602.5
$ws.Range("K132").Value = 5791.2219
$ws.Range("L132").Value = 10807.5
$ws.Range("M132").Value = -3261.2219
$ws.Range("N132").Value = -15867.5

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 65695
$ws.Range("J137").Value = 65695
$ws.Range("L137").Value = 65695
$ws.Range("N137").Value = -75895

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 267.33334
$ws.Range("I22").Value = 267.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 267.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 82.66665999999998
$ws.Range("N22").ClearContents()

$ws.Range("H58").Value = 1168.6586
$ws.Range("I58").Value = 629.88464
$ws.Range("J58").Value = 2102.5334
$ws.Range("K58").Value = 629.88464
$ws.Range("L58").Value = 2102.5334
$ws.Range("M58").Value = -426.88464
$ws.Range("N58").Value = -2508.5334

$ws.Range("H132").Value = 2670
$ws.Range("I132").Value = 1651
$ws.Range("J132").Value = 4402.3
$ws.Range("K132").Value = 4953
$ws.Range("L132").Value = 13206.9
$ws.Range("M132").Value = -2423
$ws.Range("N132").Value = -18266.9

$ws.Range("H136").Value = 1168.6586
$ws.Range("I136").Value = 629.88464
$ws.Range("J136").Value = 2102.5334
$ws.Range("K136").Value = 1889.65392
$ws.Range("L136").Value = 6307.600199999999
$ws.Range("M136").Value = 660.34608
$ws.Range("N136").Value = -11407.6002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1064.8182
$ws.Range("J5").Value = 1256.7826
$ws.Range("L5").Value = 3770.3478
$ws.Range("N5").Value = -3994.3478

$ws.Range("H8").Value = 386.375
$ws.Range("I8").Value = 386.375
$ws.Range("K8").Value = 1159.125
$ws.Range("M8").Value = -1020.125

$ws.Range("H68").Value = 1734.6383
$ws.Range("I68").Value = 1618.7812
$ws.Range("J68").Value = 1981.8
$ws.Range("K68").Value = 4856.3436
$ws.Range("L68").Value = 5945.4
$ws.Range("M68").Value = -4045.3436
$ws.Range("N68").Value = -7567.4

$ws.Range("H71").Value = 1734.6383
$ws.Range("I71").Value = 1618.7812
$ws.Range("J71").Value = 1981.8
$ws.Range("K71").Value = 14569.0308
$ws.Range("L71").Value = 17836.2
$ws.Range("M71").Value = -10513.0308
$ws.Range("N71").Value = -25948.2

$ws.Range("H98").Value = 300.42856
$ws.Range("I98").Value = 316.5
$ws.Range("K98").Value = 949.5
$ws.Range("M98").Value = 548.5

$ws.Range("H135").Value = 1064.8182
$ws.Range("J135").Value = 1256.7826
$ws.Range("L135").Value = 11311.0434
$ws.Range("N135").Value = -16381.0434

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H55").Value = 8000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 8000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8654
$ws.Range("M55").ClearContents()

$ws.Range("H139").Value = 59799.2
$ws.Range("J139").Value = 59799.2
$ws.Range("L139").Value = 59799.2
$ws.Range("N139").Value = -70079.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 100000
$ws.Range("J75").Value = 100000
$ws.Range("L75").Value = 100000
$ws.Range("N75").Value = -101872

$ws.Range("H78").Value = 100000
$ws.Range("J78").Value = 100000
$ws.Range("L78").Value = 300000
$ws.Range("N78").Value = -309360

$ws.Range("H122").Value = 3681.818
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -16150

$ws.Range("H132").Value = 4545.8
$ws.Range("I132").Value = 4480.6
$ws.Range("J132").Value = 4611
$ws.Range("K132").Value = 13441.8
$ws.Range("L132").Value = 13833
$ws.Range("M132").Value = -10911.8
$ws.Range("N132").Value = -18893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 734.4
$ws.Range("I107").Value = 693.7778
$ws.Range("K107").Value = 2081.3334
$ws.Range("M107").Value = -161.3334
Write-Host "Edit complete"
